$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename header row: "<Column>_old" -> "<Column>_FV2404" and
#    "<Column>_new" -> "<Column>_FV2410" (the "diff" header stays as-is).
# ---------------------------------------------------------------------------
$newHeaders = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $newHeaders.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newHeaders[$i]
}

# ---------------------------------------------------------------------------
# 2) Turn the A1:U77 range into a proper Excel Table (ListObject) with an
#    autofilter, keeping the header row's pre-existing look (bold font, grey
#    fill, thin borders, centered/wrapped text) unchanged and without Excel
#    auto-generating a header-row style override (dxf). We do this by
#    stashing the current header formatting, clearing it (so the table is
#    created against "default" formatting and no dxf is recorded), creating
#    the table, and then restoring the original formatting via copy/paste of
#    formats only.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$usedRange = $ws.Range("A1:U77")

# stash current header formatting using a scratch row far below the data
$scratch = $ws.Range("A200:U200")
$headerRange.Copy()
$scratch.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$headerRange.ClearFormats()

$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $usedRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# restore the header formatting exactly as it was before
$scratch.Copy()
$headerRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# clean up the scratch row completely
$scratch.ClearContents()
$scratch.ClearFormats()

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split below row 1).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
